$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "movies_data"

$ws1.Range("A2:A129").NumberFormat = "MM-DD-YYYY"
$ws2.Range("A2:A129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("A2:A129").ClearFormats()
$ws3.Range("A2:A129").NumberFormat = "mm-dd-yyyy"
